$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (new field names)
$ws.Range("A1").Value = "Distribution"
$ws.Range("B1").Value = "Specimen"
$ws.Range("C1").Value = "Bewertung"
$ws.Range("D1").Value = "Bemerkungen"
$ws.Range("E1").Value = "Massnahmen"

# Clear the example/sample data row (row 2), keep formatting
$ws.Range("A2:E2").ClearContents()

# Update the selected/active cell to D17
$ws.Range("D17").Select()
